$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (MDD 0)
$ws.Range("C2").Value = $true
$ws.Range("D2").Value = 0.06380320388317115
$ws.Range("E2").Value = 0.06380320388317115

# Row 4 (MDD 46)
$ws.Range("D4").Value = 0.2543534156797227
$ws.Range("E4").Value = 0.2543534156797227

# Row 5 (MDD 17)
$ws.Range("D5").Value = [double]"2.107138718443732E-101"
$ws.Range("E5").Value = [double]"2.107138718443732E-101"

# Row 6 (MDD 23)
$ws.Range("D6").Value = 0.0004047330648454061
$ws.Range("E6").Value = 0.0004047330648454061

# Row 8 (Control 13)
$ws.Range("D8").Value = 0.9999999999999973
$ws.Range("E8").Value = [double]"2.664535259100376E-15"

# Row 10 (Control 51)
$ws.Range("D10").Value = [double]"2.142946429496777E-16"
$ws.Range("E10").Value = 0.9999999999999998

# Row 11 (Success %)
$ws.Range("D11").Value = 0.06891186566726608
$ws.Range("E11").Value = 0.931088134332734
$ws.Range("F11").Value = 183.2982330322266
$ws.Range("G11").Value = 0.6
